{"js": "// Update the CV (\"Se actualiza la CV\"):\n//  1. Replace the old GitHub Pages portfolio link with the new Netlify one.\n//  2. Drop \"pr\u00e1cticas o \" from the \"primer empleo\" bullet under\n//     \"Otros datos de inter\u00e9s\".\n\nconst body = context.document.body;\n\n// --- 1) Portfolio link -----------------------------------------------\nconst oldLink = \"https://kevolive.github.io/Portafolio_KAOF/\";\nconst newLink = \"https://portafoliokevinolivella.netlify.app/#inicio\";\n\nconst linkResults = body.search(oldLink, { matchCase: true, matchWholeWord: false });\nlinkResults.load(\"text\");\nawait context.sync();\n\nif (linkResults.items.length === 0) {\n  throw new Error(\"Could not find the old portfolio link text to replace.\");\n}\n\nfor (let i = 0; i < linkResults.items.length; i++) {\n  linkResults.items[i].insertText(newLink, Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2) \"Inter\u00e9s en pr\u00e1cticas o primer empleo...\" bullet -------------\nconst oldBullet = \"Inter\u00e9s en pr\u00e1cticas o primer empleo en desarrollo web\";\nconst newBullet = \"Inter\u00e9s en primer empleo en desarrollo web\";\n\nconst bulletResults = body.search(oldBullet, { matchCase: true, matchWholeWord: false });\nbulletResults.load(\"text\");\nawait context.sync();\n\nif (bulletResults.items.length === 0) {\n  throw new Error(\"Could not find the bullet text to update.\");\n}\n\nfor (let i = 0; i < bulletResults.items.length; i++) {\n  bulletResults.items[i].insertText(newBullet, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Update the CV (\"Se actualiza la CV\"):\n#  1. Replace the old GitHub Pages portfolio link with the new Netlify one.\n#  2. Drop \"pr\u00e1cticas o \" from the \"primer empleo\" bullet under\n#     \"Otros datos de inter\u00e9s\".\n\n$d = $word.ActiveDocument\n\n# --- 1) Portfolio link -------------------------------------------------\n$oldLink = \"https://kevolive.github.io/Portafolio_KAOF/\"\n$newLink = \"https://portafoliokevinolivella.netlify.app/#inicio\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$found = $find.Execute($oldLink, $false, $false, $false, $false, $false, $true, 1, $false, $newLink, 2)\nif (-not $found) {\n    throw \"Could not find the old portfolio link text to replace.\"\n}\n\n# --- 2) \"Inter\u00e9s en pr\u00e1cticas o primer empleo...\" bullet ---------------\n$oldBullet = \"Inter\u00e9s en pr\u00e1cticas o primer empleo en desarrollo web\"\n$newBullet = \"Inter\u00e9s en primer empleo en desarrollo web\"\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$found2 = $find2.Execute($oldBullet, $false, $false, $false, $false, $false, $true, 1, $false, $newBullet, 2)\nif (-not $found2) {\n    throw \"Could not find the bullet text to update.\"\n}\n"}
